$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 5.112200000000002
$ws.Range("B6").Value = 9.374100000000002
$ws.Range("B7").Value = 6.249499999999997
$ws.Range("B16").Value = 8.993500000000008
$ws.Range("B20").Value = 5.841399999999997
